$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# List of (cell, new text value) updates coming from the upstream
# "Updated symbol list" GitHub Actions data refresh. Every cell in this
# sheet is stored as plain text (coin price / % change strings), so we
# force the target range to Text format ("@") before writing the new
# value -- this prevents Excel from "helpfully" re-interpreting a
# numeric-looking string (e.g. "287.89") or a percent-looking string
# (e.g. "-0.81%") as a Number/Percentage cell, which would silently
# change the cell type and strip significant trailing zeros.
$updates = @(
    @{ Cell = "D2"; Value = "287.89" }
    @{ Cell = "E2"; Value = "-0.81%" }
    @{ Cell = "D3"; Value = "30.97" }
    @{ Cell = "E3"; Value = "1.55%" }
    @{ Cell = "D4"; Value = "4.913" }
    @{ Cell = "E4"; Value = "-0.78%" }
    @{ Cell = "D5"; Value = "0.07319" }
    @{ Cell = "E5"; Value = "1.49%" }
    @{ Cell = "D6"; Value = "2.341" }
    @{ Cell = "E6"; Value = "27.40%" }
    @{ Cell = "D7"; Value = "7.724" }
    @{ Cell = "E7"; Value = "0.44%" }
    @{ Cell = "D8"; Value = "3.723" }
    @{ Cell = "E8"; Value = "-1.22%" }
    @{ Cell = "D9"; Value = "0.9031" }
    @{ Cell = "E9"; Value = "0.71%" }
    @{ Cell = "D10"; Value = "0.09126" }
    @{ Cell = "E10"; Value = "17.97%" }
    @{ Cell = "D11"; Value = "0.1692" }
    @{ Cell = "E11"; Value = "1.74%" }
    @{ Cell = "D12"; Value = "0.08153" }
    @{ Cell = "E12"; Value = "2.14%" }
    @{ Cell = "D13"; Value = "0.03127" }
    @{ Cell = "E13"; Value = "2.93%" }
    @{ Cell = "D14"; Value = "0.09936" }
    @{ Cell = "E14"; Value = "-0.81%" }
    @{ Cell = "D15"; Value = "0.001496" }
    @{ Cell = "E15"; Value = "0.30%" }
    @{ Cell = "D16"; Value = "0.005822" }
    @{ Cell = "E16"; Value = "-0.59%" }
    @{ Cell = "D17"; Value = "3.493" }
    @{ Cell = "E17"; Value = "0.79%" }
    @{ Cell = "D18"; Value = "2.099" }
    @{ Cell = "E18"; Value = "0.88%" }
    @{ Cell = "E19"; Value = "0.31%" }
    @{ Cell = "D20"; Value = "0.1293" }
    @{ Cell = "E20"; Value = "1.11%" }
    @{ Cell = "D21"; Value = "4.195" }
    @{ Cell = "E21"; Value = "3.79%" }
    @{ Cell = "E22"; Value = "-12.15%" }
    @{ Cell = "D23"; Value = "0.04511" }
    @{ Cell = "E23"; Value = "0.09%" }
    @{ Cell = "D24"; Value = "0.001212" }
    @{ Cell = "E24"; Value = "-0.42%" }
    @{ Cell = "E25"; Value = "-9.98%" }
    @{ Cell = "E26"; Value = "3.83%" }
    @{ Cell = "D27"; Value = "0.0003395" }
    @{ Cell = "D39"; Value = "0.01576" }
    @{ Cell = "E39"; Value = "0.51%" }
    @{ Cell = "D40"; Value = "0.04444" }
    @{ Cell = "E40"; Value = "1.50%" }
    @{ Cell = "D41"; Value = "0.007335" }
    @{ Cell = "E41"; Value = "0.06%" }
    @{ Cell = "D42"; Value = "0.009502" }
    @{ Cell = "E42"; Value = "-4.06%" }
    @{ Cell = "E43"; Value = "1.99%" }
    @{ Cell = "E44"; Value = "10.63%" }
    @{ Cell = "D45"; Value = "0.009001" }
    @{ Cell = "E45"; Value = "-5.59%" }
    @{ Cell = "D46"; Value = "0.00006102" }
    @{ Cell = "E46"; Value = "2.34%" }
    @{ Cell = "E47"; Value = "-0.15%" }
    @{ Cell = "D48"; Value = "2.330" }
    @{ Cell = "E48"; Value = "0.88%" }
    @{ Cell = "E50"; Value = "-0.15%" }
    @{ Cell = "E51"; Value = "-0.15%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}

